$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 7: shift data (was Dina/Primera @44161) to Castle Brite/Especial @44546
$ws.Range("D7").Value = 44546
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("L7").Value = "Especial"
$ws.Range("N7").Value = 22500
$ws.Range("O7").Value = 23000
$ws.Range("P7").Value = 22750
$ws.Range("Q7").Value = "$/caja 18 kilos"
$ws.Range("S7").Value = 1264
$ws.Range("T7").Value = 18

# Update row 8
$ws.Range("D8").Value = 44546
$ws.Range("K8").Value = "Castle Brite"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 20500
$ws.Range("O8").Value = 21000
$ws.Range("P8").Value = 20750
$ws.Range("Q8").Value = "$/caja 18 kilos"
$ws.Range("S8").Value = 1153
$ws.Range("T8").Value = 18

# Update row 9
$ws.Range("D9").Value = 44161
$ws.Range("K9").Value = "Dina"
$ws.Range("M9").Value = 300
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20500
$ws.Range("P9").Value = 20250
$ws.Range("S9").Value = 1350

# Update row 10
$ws.Range("D10").Value = 44161
$ws.Range("K10").Value = "Dina"
$ws.Range("L10").Value = "Segunda"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18500
$ws.Range("P10").Value = 18250
$ws.Range("Q10").Value = "$/caja 15 kilos"
$ws.Range("S10").Value = 1217
$ws.Range("T10").Value = 15

# New row 11
$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C11").Value = "Coquimbo"
$ws.Range("D11").Value = 44160
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100103
$ws.Range("H11").Value = "Frutos de hueso (carozo)"
$ws.Range("I11").Value = 100103003
$ws.Range("J11").Value = "Damasco"
$ws.Range("K11").Value = "Castle Brite"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 240
$ws.Range("N11").Value = 20500
$ws.Range("O11").Value = 21000
$ws.Range("P11").Value = 20750
$ws.Range("Q11").Value = "$/caja 15 kilos"
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 1383
$ws.Range("T11").Value = 15

# New row 12
$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C12").Value = "Coquimbo"
$ws.Range("D12").Value = 44175
$ws.Range("E12").Value = 4
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100103
$ws.Range("H12").Value = "Frutos de hueso (carozo)"
$ws.Range("I12").Value = 100103003
$ws.Range("J12").Value = "Damasco"
$ws.Range("K12").Value = "Castle Brite"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 300
$ws.Range("N12").Value = 21000
$ws.Range("O12").Value = 22000
$ws.Range("P12").Value = 21500
$ws.Range("Q12").Value = "$/caja 18 kilos"
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 1194
$ws.Range("T12").Value = 18

# Copy the date style from D10 (existing) to D11:D12 so they keep the same number format
$ws.Range("D9").Copy()
$ws.Range("D11:D12").PasteSpecial(-4122) | Out-Null
